$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new chef item (Test Chef)
$ws.Range("A3").Value = 13
$ws.Range("B3").Value = "testchef@gmail.com"
$ws.Range("C3").Value = "Test"
$ws.Range("D3").Value = "Chef"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = "Chef"
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = $false

# Row 4: new chef item (chef test)
$ws.Range("A4").Value = 15
$ws.Range("B4").Value = "chef@gmail.com"
$ws.Range("C4").Value = "chef"
$ws.Range("D4").Value = "test"
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = "Chef"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 20
$ws.Range("I4").Value = $true

$ws.Columns("A:H").AutoFit() | Out-Null
